$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the two shapes that reference the Cummings & Theodore (2022) citation:
#  - the in-line body-text mention: "...(Cummings & Theodore, 2022)."
#  - the full reference-list entry: "Cummings, S. N. & Theodore, R. M. (2022). ..."
$bodyShape = $null
$refShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if (-not $shp.HasTextFrame) { continue }
    if (-not $shp.TextFrame.HasText) { continue }
    $t = $shp.TextFrame.TextRange.Text
    if ($t.Contains("(Cummings &")) {
        $bodyShape = $shp
    } elseif ($t.Contains("Cummings, S. N.")) {
        $refShape = $shp
    }
}

# --- In-line citation: "(Cummings & Theodore, 2022)." -> "(Cummings & Theodore, accepted)."
$bodyTr = $bodyShape.TextFrame.TextRange
$bodyFull = $bodyTr.Text

# Remove the "2022" year, leaving "(Cummings & Theodore, )."
$yearIdx = $bodyFull.IndexOf("2022")
$yearRange = $bodyTr.Characters($yearIdx + 1, 4)
$yearRange.Text = ""

# Re-read text, then turn the remaining (italic) "." run into "accepted)."
$bodyFull2 = $bodyTr.Text
$closeIdx = $bodyFull2.IndexOf(").", $yearIdx)
$closeRange = $bodyTr.Characters($closeIdx + 1, 2)
$closeRange.Text = "accepted)."

# --- Reference-list entry: "(2022)." -> "(accepted)." (kept italic, matching the journal-name run style)
$refTr = $refShape.TextFrame.TextRange
$refFull = $refTr.Text

$refIdx = $refFull.IndexOf("(2022)")
$refRange = $refTr.Characters($refIdx + 1, 8)
$refRange.Font.Italic = $true
$refRange.Text = "(accepted). "
